# Generate Report for Handoff
# - Update status text from "In Translation" to "Ready for handoff"
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Widen the status-related columns to fit the new, longer status text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Update status values on the Overview sheet (zh-cn / de-de status columns) ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-02 11:09:10"

# --- Update status + handoff datetime on the per-language sheets ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-02 11:09:00"

$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-02 11:09:10"

# --- Widen the columns that now hold the longer "Ready for handoff" text ---
# (ColumnWidth is quantized to whole pixels by the Excel engine, so the input
# value below is the closest representable width to the target 17.2159881591797
# characters; it rounds internally to the nearest achievable column width.)
$newColWidth = 16.333333333333332
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

$zhcn.Columns.Item(3).ColumnWidth = $newColWidth
$dede.Columns.Item(3).ColumnWidth = $newColWidth
